$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 5 (a2, a3, and the old "my a" rows)
$ws.Range("A3:C5").EntireRow.Delete()

# Update row 2 to hold the "my a" record with the new C value
$ws.Range("A2").Value = "my a"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 7
